$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 97) with the latest OHLC data, mirroring the
# existing rows produced by the R script that feeds this sheet.
$row = 97
$prev = $row - 1

# Column A: date/time serial value. Copy the previous row's cell formatting
# (date/time number format) so the new cell reuses the same style instead of
# minting a new one.
$ws.Cells.Item($row, 1).Value = 45462.2916666667
$ws.Cells.Item($prev, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats

# Column B: volume
$ws.Cells.Item($row, 2).Value = 0

# Columns C-F: high, low, open, close (numeric)
$ws.Cells.Item($row, 3).Value = 0.699999988079071
$ws.Cells.Item($row, 4).Value = 0.699999988079071
$ws.Cells.Item($row, 5).Value = 0.699999988079071
$ws.Cells.Item($row, 6).Value = 0.699999988079071

# Column G: adj_close is stored as text in this sheet (matches the source
# data, which came straight out of an R script rather than Excel). Force the
# text number format so the numeric-looking literal isn't auto-coerced into
# a number.
$ws.Cells.Item($row, 7).NumberFormat = "@"
$ws.Cells.Item($row, 7).Value = "0.699999988079071"

# Column H: ticker (text)
$ws.Cells.Item($row, 8).Value = "BWZ.MI"

$excel.CutCopyMode = $false
